# Commit: "Updated symbol list on Mon Feb  6 07:51:00 UTC 2023 with GitHub Actions"
#
# This script refreshes the live crypto-price table on Sheet1:
#   - Columns B (Coin) and C (Link) are updated for rows 18-25, where the
#     site's coin ranking reshuffled: the former row 25 entry (HotbitToken)
#     moved up to row 18, and every entry from the old rows 18-24 shifted
#     down by one row.
#   - Column D (Price) and column E (Volume 1h) are refreshed with the
#     latest scraped values for many rows across the sheet.
#
# All of these cells hold text (not real numbers) in the workbook - e.g.
# "0.08060" (with a significant trailing zero) and "-2.16%" - so plain
# numeric assignment must be avoided (Excel would normalize "0.08060" to
# 0.0806, or turn "-2.16%" into a percentage-formatted number). Each
# numeric-looking cell is therefore temporarily switched to the Text
# number format before the value is written, and its original style is
# restored immediately afterwards so no formatting is altered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $savedStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = $savedStyle
}

# --- Coin / Link reshuffle (rows 18-25) ---
# Plain text, not numeric-looking, so a direct assignment is fine.
$ws.Range('B18').Value = 'HotbitToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('B19').Value = 'LEO'
$ws.Range('C19').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('B20').Value = 'BitpandaEcosystemToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('B21').Value = 'MCDex'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('B22').Value = 'ProBitToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('B23').Value = 'ZBToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('B24').Value = 'CoinExToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('B25').Value = 'BitKan'
$ws.Range('C25').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'

# --- Price (D) / Volume 1h (E) refresh ---
$priceVolumeUpdates = [ordered]@{
    'D2'  = '325.31';   'E2'  = '-2.16%'
    'D3'  = '44.47';    'E3'  = '0.86%'
    'D4'  = '5.494';    'E4'  = '-6.22%'
    'D5'  = '0.08060';  'E5'  = '-3.28%'
    'D6'  = '8.640';    'E6'  = '-1.99%'
    'D7'  = '1.905';    'E7'  = '-3.32%'
    'D8'  = '4.276';    'E8'  = '-4.98%'
    'D10' = '0.9409';   'E10' = '0.54%'
    'D11' = '0.1171';   'E11' = '-6.95%'
    'D12' = '0.1867';   'E12' = '-4.67%'
    'D13' = '0.1005';   'E13' = '4.23%'
    'D14' = '0.04258';  'E14' = '3.63%'
    'D15' = '0.1064';   'E15' = '-0.17%'
    'D16' = '0.001278'; 'E16' = '-2.62%'
    'D17' = '0.005871'; 'E17' = '-0.75%'
    'D18' = '0.004556'; 'E18' = '3.27%'
    'D19' = '3.584';    'E19' = '2.27%'
    'D20' = '0.3498';   'E20' = '-0.34%'
    'D21' = '8.445';    'E21' = '-3.86%'
    'D22' = '0.1378';   'E22' = '0.57%'
    'D23' = '0.2528';   'E23' = '-1.64%'
    'D24' = '0.04236';  'E24' = '-3.74%'
    'D25' = '0.001235'; 'E25' = '-1.81%'
    'E26' = '-0.94%'
    'D27' = '0.0003990'; 'E27' = '-0.01%'
    'E39' = '-6.51%'
    'D40' = '0.05478';   'E40' = '-3.94%'
    'D41' = '0.007670';  'E41' = '-3.07%'
    'E42' = '-2.22%'
    'D43' = '0.007156';  'E43' = '-20.81%'
    'D44' = '0.002016';  'E44' = '-4.12%'
    'D45' = '0.009197';  'E45' = '-12.26%'
    'D46' = '0.00007093'; 'E46' = '-2.21%'
    'E47' = '-0.01%'
    'D48' = '0.003586';  'E48' = '10.87%'
    'D49' = '0.002271';  'E49' = '-0.41%'
    'D50' = '0.00002101'; 'E50' = '-0.01%'
    'D51' = '0.0002001'; 'E51' = '-0.01%'
}

foreach ($cellRef in $priceVolumeUpdates.Keys) {
    Set-TextValue $cellRef $priceVolumeUpdates[$cellRef]
}
